$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 (pushes the old merged "invalid row" message
# down to row 7); Excel carries the row-5 formatting down into the new row 6.
$ws.Rows.Item(6).Insert()

# The new row's description cell (C) wraps instead of shrink-to-fit, and the
# row is taller to show the wrapped text.
$ws.Range("C6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 30

# New sample data row.
$ws.Range("A6").Value = "A005"
$ws.Range("B6").Value = 'TOMES_PATTERN: {"A","B"}, {1,"2"}'
$ws.Range("C6").Value = 'Should yield 4 (2*2) manifestations AFTER first bing unable to sort due to the "1" not being in quotes.'
$ws.Range("D6").Value = $true
$ws.Range("E6").Value = "TOMES"
$ws.Range("F6").Value = "ncdcr.gov"

# Give the new B6 cell the same "duplicate value" / "contains [" / "contains
# .txt" conditional-formatting rules the other B-column data cells have.
$fcs = $ws.Range("B6").FormatConditions

$fcDup = $fcs.AddUniqueValues()
$fcDup.DupeUnique = 1
$fcDup.Font.Color = 22428
$fcDup.Interior.Color = 10284031

$fcBracket = $fcs.Add(9)
$fcBracket.Modify(9, 0, $null, "[")
$fcBracket.Text = "["
$fcBracket.Formula1 = '=NOT(ISERROR(SEARCH("[",B6)))'
$fcBracket.Font.Color = 24832
$fcBracket.Interior.Color = 13561798

$fcTxt = $fcs.Add(9)
$fcTxt.Modify(9, 0, $null, ".txt")
$fcTxt.Text = ".txt"
$fcTxt.Formula1 = '=NOT(ISERROR(SEARCH(".txt",B6)))'
$fcTxt.Font.Color = 393372
$fcTxt.Interior.Color = 13551615

# Re-establish rule priority so the new rules sit at the top (1,2,3), matching
# how Excel bumps existing rules down when new ones are added.
$fcTxt.SetFirstPriority()
$fcBracket.SetFirstPriority()
$fcDup.SetFirstPriority()

# Move the active selection to the cell the author ended up editing last.
$ws.Range("D6").Select()
